# 自动更新价格数据: insert today's price row at the top of the table,
# pushing all existing date rows down by one (latest date first).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row just below the header row (row 1), shifting the
# existing data (and all the dates below it) down by one row.
$ws.Rows.Item(2).Insert()

# Populate the new row with the latest date and the same metric values
# that used to sit in row 2 (prices unchanged day-over-day).
# Use a leading apostrophe so Excel stores the date as literal text
# rather than converting it to a date serial number.
$ws.Range("A2").Value = "'2026-01-05"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# The freshly inserted row inherits formatting copied down from the
# header row (bold/centered/bordered). Clear it so the new data row
# matches the plain, unstyled look of every other data row.
$ws.Range("A2:D2").ClearFormats()
